$d = $word.ActiveDocument

# --- Paragraph: "Horas del día 4: 2h" ---
# InsertParagraphAfter() on the last paragraph creates a new paragraph that
# inherits the pPr/rPr (jc=both, rFonts cstheme=minorHAnsi) of the paragraph
# it follows, and carries that formatting onto the run created when text is
# assigned to it.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter() | Out-Null
$horasPara = $d.Paragraphs.Last
$horasPara.Range.Text = "Horas del día 4: 2h"

# --- Blank separator paragraph (no run at all) ---
# Using Find/Replace with the "^p" special sequence to append a paragraph
# mark (instead of Range.InsertParagraphAfter) produces a genuinely empty
# paragraph with no <w:r> child, matching the other blank separator
# paragraphs already present in this document.
$d.Content.Find.Execute("Horas del día 4: 2h", $true, $false, $false, $false,
    $false, $true, 1, $false, "Horas del día 4: 2h^p", 2) | Out-Null

# --- Paragraph: "-Día 5: ..." (three runs) ---
$dia5Para = $d.Paragraphs.Last
$dia5Para.Range.InsertParagraphAfter() | Out-Null
$dia5Para = $d.Paragraphs.Last

$t1 = "-Día 5: encuentro la manera de solucionar las funciones copiar y cortar, utilizando QtGui.QTextCursor"
$t2 = ". Añado la funcionalidad pegar"
$t3 = ". He completado mi código con un visualizador de markdown en formato HTML, teniendo ahora dos aplicaciones. Mi próximo paso es añadir la app de visualización de markdown a mi aplicación de edición."

# Assign the full concatenated text first so the whole paragraph gets a
# single, uniformly-formatted run inheriting rFonts cstheme=minorHAnsi.
$dia5Para.Range.Text = $t1 + $t2 + $t3

# Now split that single run into three separate <w:r> elements (matching
# the target, which has three runs with identical rPr) by dropping a
# bookmark at each boundary and immediately deleting it: inserting /
# removing a bookmark forces a run split without altering the run
# formatting or leaving any bookmark markup behind.
$paraStart = $dia5Para.Range.Start
$split1 = $paraStart + $t1.Length
$split2 = $split1 + $t2.Length

$d.Bookmarks.Add("zzSplit1", $d.Range($split1, $split1)) | Out-Null
$d.Bookmarks.Add("zzSplit2", $d.Range($split2, $split2)) | Out-Null
$d.Bookmarks("zzSplit1").Delete() | Out-Null
$d.Bookmarks("zzSplit2").Delete() | Out-Null
